# Updated Karma Specs, tests written
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Update the three "X" cells that should now read "Y-Test"
$ws.Range("F12").Value = "Y-Test"
$ws.Range("C13").Value = "Y-Test"
$ws.Range("F13").Value = "Y-Test"
$ws.Range("C14").Value = "Y-Test"

# Move the active cell / selection to F14
$ws.Range("F14").Select()
